# The "last updated" date column (C) for every data row (rows 2-533)
# is bumped by one day: 45178 (2023-09-09) -> 45179 (2023-09-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C533").Value = 45179
